# "Generate Report for Handback"
#
# The handback transform for the c5067ff2-... file failed because the
# handback file name did not match the handoff file name. Record that
# failure on the Overview sheet (status text) and write the detailed
# error message into the "Error Detail" column of the per-locale
# (zh-cn / de-de) report rows, widening that column so the message is
# readable.

$wb = $excel.ActiveWorkbook

$zhMessage = "Handback file name: fjjn5dav.who is different with handoff file name: c5067ff2-1fbb-411e-8cd5-1e4e87c0ddd7.ff7bc03cb9c9d08b97d8c617079c6792c83f879c.zh-cn."
$deMessage = "Handback file name: fjjn5dav.who is different with handoff file name: c5067ff2-1fbb-411e-8cd5-1e4e87c0ddd7.ff7bc03cb9c9d08b97d8c617079c6792c83f879c.de-de."

$statusText = "Handback transform failed"

# Overview sheet: the c5067ff2 row's status (shared by zh-cn & de-de
# columns) moves from "Ready for handoff" to "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# The saved OOXML column width ends up offset from COM's ColumnWidth by
# 5/6 (the default-font padding) in this runtime, so asking for
# 39.1666... here is what lands on a clean "width=40" in the XML.
$errorColumnWidth = 39.1666666666667

# zh-cn sheet: row 3 (c5067ff2 file) - same Status text, plus the new
# Error Detail column (P) value.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhMessage
$wsZh.Columns.Item(16).ColumnWidth = $errorColumnWidth

# de-de sheet: row 3 (c5067ff2 file) - same Status text, plus the new
# Error Detail column (P) value.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deMessage
$wsDe.Columns.Item(16).ColumnWidth = $errorColumnWidth
